$d = $word.ActiveDocument

$replacements = @(
    @("523÷8=", "327÷5="),
    @("571÷3=", "886÷4="),
    @("262÷5=", "539÷3="),
    @("845÷6=", "271÷9="),
    @("362÷8=", "283÷5="),
    @("733÷2=", "920÷4="),
    @("991÷9=", "745÷8="),
    @("988÷9=", "943÷9="),
    @("746÷4=", "369÷6="),
    @("431÷2=", "733÷8="),
    @("200÷9=", "623÷3="),
    @("172÷5=", "628÷8="),
    @("762÷5=", "745÷7="),
    @("521÷4=", "662÷6="),
    @("528÷8=", "447÷8="),
    @("996÷3=", "910÷2="),
    @("395÷3=", "817÷8="),
    @("189÷6=", "492÷9="),
    @("128÷9=", "627÷9="),
    @("539÷2=", "649÷5="),
    @("270÷4=", "399÷8="),
    @("109÷8=", "803÷2="),
    @("475÷9=", "351÷8="),
    @("287÷5=", "665÷2="),
    @("669÷5=", "611÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
